$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.226658225059509
$ws.Range("B1").Value = 2.778029680252075
$ws.Range("C1").Value = 4.756212711334229
$ws.Range("D1").Value = 2.09923529624939
$ws.Range("E1").Value = 1.157722353935242
